$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 3, shifting rows 3-5 down to 4-6
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new descriptive strings
$ws.Range("A3").Value = "User updated record"
$ws.Range("B3").Value = "Timestamp of record update"
